# Add a new test-case row (row 20) to the "Test Cases" sheet for DRA0017,
# mirroring the layout/formatting of the row directly above it (row 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values -----------------------------------------------------------
$ws.Range("A20").Value = "DRA0017"
$ws.Range("B20").Value = "OPQA-4525||OPQA-4526||OPQA-4527"
$ws.Range("C20").Value = "Verify that the STeAM Step Up Auth Modal should be presented to the user without a pre-populated email address when user has a valid Neon session token and is navigating within the same browser window."
$ws.Range("D20").Value = "Y"

# --- Formatting ---------------------------------------------------------
# Clone the border/fill/number-format of the row above (A,B,D,E columns
# share the plain bordered style; C has the wrapped/bordered style) so the
# new row visually matches the rest of the table.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)

$ws.Range("B19").Copy()
$ws.Range("B20").PasteSpecial(-4122)

$ws.Range("C19").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("D19").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E19").Copy()
$ws.Range("E20").PasteSpecial(-4122)

# Description cell gets its own (slightly darker) font color and keeps
# the wrap so the long text is readable. (Only touch .Color - touching
# .Name/.Size here would force a full font rewrite that drops the
# inherited "minor scheme" flag.)
$descCell = $ws.Range("C20")
$descCell.Font.Color = 3355443
$descCell.WrapText = $true

# Row height to fit the wrapped description (3 lines).
$ws.Rows.Item(20).RowHeight = 45

# --- View state (best effort) -------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("C28").Select()
